# Adding PC SunEnergy to the Portfolio Forecast
# Update the Consumption Forecast sheet: new forecast values (column A) and
# timestamps shifted forward by 2 days (column B), for rows 2-97.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newConsumption = @(
    5330,5280,5220,5170,5130,5090,5060,5050,5050,5050,
    5050,5050,5050,5050,5060,5070,5100,5140,5190,5260,
    5350,5440,5540,5640,5730,5820,5900,5950,5990,6010,
    6010,5980,5940,5880,5820,5750,5680,5610,5550,5500,
    5460,5420,5390,5370,5340,5310,5290,5270,5250,5240,
    5240,5250,5270,5300,5320,5350,5380,5410,5440,5470,
    5510,5560,5620,5680,5750,5830,5900,5970,6040,6120,
    6200,6300,6400,6500,6600,6700,6790,6880,6960,7020,
    7060,7070,7030,6920,6800,6650,6480,6330,6150,6000,
    5860,5720,5460,5400,5350,5300
)

$dayShift = 2

for ($i = 0; $i -lt $newConsumption.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $newConsumption[$i]
    $ws.Cells.Item($row, 2).Value = $ws.Cells.Item($row, 2).Value2 + $dayShift
}
